$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '37.090.64'
$c.ClearFormats()
$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '
$c.ClearFormats()

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.051.49'
$c.ClearFormats()
$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  -0.51%  '
$c.ClearFormats()

$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.ClearFormats()

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '249.44'
$c.ClearFormats()
$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  -0.09%  '
$c.ClearFormats()

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.669'
$c.ClearFormats()
$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  -0.78%  '
$c.ClearFormats()

$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '59.25'
$c.ClearFormats()
$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  +8.39%  '
$c.ClearFormats()

$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  +0.00%  '
$c.ClearFormats()

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.388'
$c.ClearFormats()
$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  +1.01%  '
$c.ClearFormats()

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.0796'
$c.ClearFormats()
$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  -0.70%  '
$c.ClearFormats()

$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  +1.87%  '
$c.ClearFormats()

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '16.04'
$c.ClearFormats()
$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  +6.63%  '
$c.ClearFormats()

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '2.347.47'
$c.ClearFormats()
$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  -0.62%  '
$c.ClearFormats()

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '0.837'
$c.ClearFormats()
$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  +2.33%  '
$c.ClearFormats()

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '5.74'
$c.ClearFormats()
$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  +7.63%  '
$c.ClearFormats()

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '2.051.42'
$c.ClearFormats()
$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  -0.54%  '
$c.ClearFormats()

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '18.40'
$c.ClearFormats()
$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  +29.04%  '
$c.ClearFormats()

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '37.073.89'
$c.ClearFormats()
$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  +0.26%  '
$c.ClearFormats()

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '75.43'
$c.ClearFormats()
$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  +2.53%  '
$c.ClearFormats()

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '0.0₃0905'
$c.ClearFormats()
$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  -3.08%  '
$c.ClearFormats()

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '5.42'
$c.ClearFormats()
$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  +0.51%  '
$c.ClearFormats()

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '238.01'
$c.ClearFormats()
$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  +0.17%  '
$c.ClearFormats()

$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c.ClearFormats()

$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '
$c.ClearFormats()

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.22'
$c.ClearFormats()
$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  +10.52%  '
$c.ClearFormats()

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '9.49'
$c.ClearFormats()
$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  +4.75%  '
$c.ClearFormats()

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '169.22'
$c.ClearFormats()
$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  -0.46%  '
$c.ClearFormats()

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '20.13'
$c.ClearFormats()
$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  +0.05%  '
$c.ClearFormats()

$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '
$c.ClearFormats()

$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  +6.28%  '
$c.ClearFormats()

$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '4.80'
$c.ClearFormats()
$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  +3.92%  '
$c.ClearFormats()

$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -0.61%  '
$c.ClearFormats()

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '4.52'
$c.ClearFormats()
$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  +2.27%  '
$c.ClearFormats()

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.0892'
$c.ClearFormats()
$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  -0.76%  '
$c.ClearFormats()

$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  -0.07%  '
$c.ClearFormats()

$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.ClearFormats()
$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  -2.33%  '
$c.ClearFormats()

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '1.75'
$c.ClearFormats()
$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -0.96%  '
$c.ClearFormats()

$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  +4.23%  '
$c.ClearFormats()

$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  -0.99%  '
$c.ClearFormats()

$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  +10.37%  '
$c.ClearFormats()

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '5.12'
$c.ClearFormats()
$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  +23.96%  '
$c.ClearFormats()

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '17.71'
$c.ClearFormats()
$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c.ClearFormats()

$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  -1.21%  '
$c.ClearFormats()

$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  -0.81%  '
$c.ClearFormats()

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '97.04'
$c.ClearFormats()
$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  -0.15%  '
$c.ClearFormats()

$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  +4.55%  '
$c.ClearFormats()

$c = $ws.Range('B47')
$c.NumberFormat = "@"
$c.Value = 'FTXToken'
$c.ClearFormats()
$c = $ws.Range('C47')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c.ClearFormats()
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '3.86'
$c.ClearFormats()
$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -7.17%  '
$c.ClearFormats()

$c = $ws.Range('B48')
$c.NumberFormat = "@"
$c.Value = 'Maker'
$c.ClearFormats()
$c = $ws.Range('C48')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c.ClearFormats()
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.289.35'
$c.ClearFormats()
$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  -0.62%  '
$c.ClearFormats()

$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  -1.21%  '
$c.ClearFormats()

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '6.81'
$c.ClearFormats()
$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  -0.99%  '
$c.ClearFormats()

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '2.223.37'
$c.ClearFormats()
$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  -1.08%  '
$c.ClearFormats()
